$d = $word.ActiveDocument

# Reference list template used by the existing numbered list (numId=1)
# e.g. paragraph "Ukol c.10-12:" already uses pStyle=Odstavecseseznamem + numPr(ilvl=0,numId=1)
$refPara = $d.Paragraphs.Item(50)
$listTemplate = $refPara.Range.ListFormat.ListTemplate

# --- Change 1: turn the "Vysledna finalni tabulka ..." paragraph into a list item ---
$p53 = $d.Paragraphs.Item(53).Range
$p53.Style = "Odstavecseseznamem"
$p53.ListFormat.ApplyListTemplate($listTemplate, $true)

# --- Change 2: remove the long "(v_base as v1, ... v_10 as v) " aside, leave two spaces ---
$p53b = $d.Paragraphs.Item(53).Range
$null = $p53b.Find.Execute(" (v_base as v1, v_34567 as v2, v_8 as v8, v_9 as v9, v_10 as v) ", `
  $true, $false, $false, $false, $false, $true, 1, $false, "  ", 2)

# --- Change 3: turn the "Vytvorila jsem tabulku ..." paragraph into a list item ---
$p55 = $d.Paragraphs.Item(55).Range
$p55.Style = "Odstavecseseznamem"
$p55.ListFormat.ApplyListTemplate($listTemplate, $true)

# --- Change 4: merge "t_alena_vitkova_projekt_SQL_" + "part" into one run ---
$p55b = $d.Paragraphs.Item(55).Range
$null = $p55b.Find.Execute("t_alena_vitkova_projekt_SQL_part", $true, $false, $false, $false, $false, `
  $true, 1, $false, "t_alena_vitkova_projekt_SQL_part", 2)

# --- Change 5: "... pro Holandsko a datumy ..." -> "... pro Holandsko s datumy ..." ---
$p55c = $d.Paragraphs.Item(55).Range
$null = $p55c.Find.Execute("Holandsko a datumy", $true, $false, $false, $false, $false, `
  $true, 1, $false, "Holandsko s datumy", 2)

Write-Output "done"
